$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 numeric values that changed
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4977036666666666
$ws.Range("N2").Value = 1.493111
$ws.Range("Q2").Value = 2.558593018785333
$ws.Range("R2").Value = 23.027337169068
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Delete row 3 entirely (this also removes the now-unused "Resolving-Mac" shared string)
$ws.Rows("3").Delete()
